# Update timestamps in the localization-status report, as if a fresh
# report had just been generated for the handback.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first row (G2).
$wsOverview.Range("G2").Value = "2016-08-18 01:00:43"

# zh-cn sheet, row 2: "Latest Handoff Datetime" (H2) and
# "Latest Handback DateTime" (K2).
$wsZhCn.Range("H2").Value = "2016-08-18 01:00:37"
$wsZhCn.Range("K2").Value = "2016-08-18 01:00:53"

# de-de sheet, row 2: "Latest Handoff Datetime" (H2, shares the same
# timestamp string as Overview!G2) and "Latest Handback DateTime" (K2).
$wsDeDe.Range("H2").Value = "2016-08-18 01:00:43"
$wsDeDe.Range("K2").Value = "2016-08-18 01:01:03"
